$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.393.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.378.57'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '551.64'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.81'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.380.45'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.34%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.34'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.52'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.808.96'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.405.98'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.383.52'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.99'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '321.33'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.58%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.15'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.90%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.45'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.91'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +9.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.496.46'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.20'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.36%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '516.98'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.68%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0901'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.150'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.02%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.53'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.84%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.88'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.17%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.379'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.53'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '146.23'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.77%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.27'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.24'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.08%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.07%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.60'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0527'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.29%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.71'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.582'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0906'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.60%  '

Write-Host "Updated cells: 141"
